# Apply the "Added placeholder statements in Water Tap" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new data row (row 37) ------------------------------------------------
$ws.Range("A37").Value = "demands\dhw"
$ws.Range("C37").Value = "WTap"
$ws.Range("B37").Value = "dhw_mfh,  dhw_mfh_1hour, dhw_mfh_CCT, dhw_mfh_v1, dhw_sfh_task44"
$ws.Range("D37").Value = "Automatic Connection Feature Completed"

# Reuse the formatting already used for similar "data" rows (row 34, style with
# green fill + thin border + vertical-center alignment) for columns A, C and D.
$ws.Range("A34").Copy()
$ws.Range("A37").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A34").Copy()
$ws.Range("C37").PasteSpecial(-4122)
$ws.Range("A34").Copy()
$ws.Range("D37").PasteSpecial(-4122)

# Column B gets the same base formatting plus left alignment + wrap text, which
# creates the new cell style used only by this cell.
$ws.Range("A34").Copy()
$ws.Range("B37").PasteSpecial(-4122)
$ws.Range("B37").HorizontalAlignment = -4131   # xlLeft
$ws.Range("B37").WrapText = $true

$excel.CutCopyMode = $false

# Taller row to fit the wrapped placeholder text.
$ws.Rows.Item(37).RowHeight = 72

# --- Minor sheet cosmetics -----------------------------------------------------
# Column B was narrowed slightly to make room for the new content.
$ws.Columns.Item(2).ColumnWidth = 22.5

# Update the visible selection to the new last cell, matching where the editor
# ended up after adding the row.
$ws.Range("J37").Select() | Out-Null
